# Auto-generated: apply row-permutation edits per target diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = 112035549
$ws.Range("B2").Value2 = 77515
$ws.Range("E2").Value2 = 6425
$ws.Range("F2").Value2 = 'Garnlav'
$ws.Range("G2").Value2 = 'Alectoria sarmentosa'
$ws.Range("H2").Value2 = '(Ach.) Ach.'
$ws.Range("Q2").Value2 = 515977.3292799139
$ws.Range("R2").Value2 = 7184566.677681392
$ws.Range("S2").Value2 = 10
$ws.Range("Z2").Value2 = '10:51'
$ws.Range("AB2").Value2 = '10:51'
$ws.Range("AH2").Value2 = 'Blåbärsgranskog'
$ws.Range("AM2").Value2 = 'Gren på levande träd'
$ws.Range("AO2").Value2 = 'Branch on living tree'

# Row 3
$ws.Range("A3").Value2 = 112038134
$ws.Range("B3").Value2 = 89405
$ws.Range("D3").Value2 = 'NT'
$ws.Range("E3").Value2 = 1202
$ws.Range("F3").Value2 = 'Ullticka'
$ws.Range("G3").Value2 = 'Phellinidium ferrugineofuscum'
$ws.Range("H3").Value2 = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("Z3").Value2 = '13:27'
$ws.Range("AB3").Value2 = '13:27'
$ws.Range("AM3").Value2 = 'Liggande död trädstam, markontakt'
$ws.Range("AO3").Value2 = 'Horizontal, dead with ground contact # Picea abies'

# Row 5
$ws.Range("A5").Value2 = 112037635
$ws.Range("B5").Value2 = 89401
$ws.Range("E5").Value2 = 1108
$ws.Range("F5").Value2 = 'Harticka'
$ws.Range("G5").Value2 = 'Pelloporus leporinus'
$ws.Range("H5").Value2 = '(Fr.) Krieglst.'
$ws.Range("Q5").Value2 = 515886.4644205247
$ws.Range("R5").Value2 = 7184225.831779522
$ws.Range("Z5").Value2 = '12:06'
$ws.Range("AB5").Value2 = '12:06'
$ws.Range("AH5").Value2 = 'Blåbärsgranskog'
$ws.Range("AJ5").Value2 = 'gran'
$ws.Range("AK5").Value2 = 'Picea abies'
$ws.Range("AM5").Value2 = 'Liggande död trädstam, markontakt'
$ws.Range("AO5").Value2 = 'Horizontal, dead with ground contact # Picea abies'

# Row 6
$ws.Range("A6").Value2 = 112038082
$ws.Range("B6").Value2 = 90087
$ws.Range("D6").Value2 = 'LC'
$ws.Range("E6").Value2 = 3298
$ws.Range("F6").Value2 = 'Trådticka'
$ws.Range("G6").Value2 = 'Climacocystis borealis'
$ws.Range("H6").Value2 = '(Fr.) Kotl. & Pouzar'
$ws.Range("Q6").Value2 = 515925.2595200292
$ws.Range("R6").Value2 = 7184319.449006356
$ws.Range("Z6").Value2 = '13:22'
$ws.Range("AB6").Value2 = '13:22'

# Row 7
$ws.Range("A7").Value2 = 112035020
$ws.Range("B7").Value2 = 89401
$ws.Range("E7").Value2 = 1108
$ws.Range("F7").Value2 = 'Harticka'
$ws.Range("G7").Value2 = 'Pelloporus leporinus'
$ws.Range("H7").Value2 = '(Fr.) Krieglst.'
$ws.Range("Q7").Value2 = 515923.0367052297
$ws.Range("R7").Value2 = 7184658.938780431
$ws.Range("S7").Value2 = 50
$ws.Range("Z7").Value2 = '10:24'
$ws.Range("AB7").Value2 = '10:24'
$ws.Range("AJ7").ClearContents()
$ws.Range("AK7").ClearContents()
$ws.Range("AM7").Value2 = 'Stubbe'
$ws.Range("AO7").Value2 = 'Stump'

# Row 8
$ws.Range("A8").Value2 = 112037386
$ws.Range("B8").Value2 = 89423
$ws.Range("E8").Value2 = 5432
$ws.Range("F8").Value2 = 'Granticka'
$ws.Range("G8").Value2 = 'Porodaedalea chrysoloma'
$ws.Range("H8").Value2 = '(Fr.) Fiasson & Niemelä'
$ws.Range("I8").ClearContents()
$ws.Range("J8").ClearContents()
$ws.Range("Q8").Value2 = 516031.6636387278
$ws.Range("R8").Value2 = 7184226.601435129
$ws.Range("Z8").Value2 = '11:52'
$ws.Range("AB8").Value2 = '11:52'
$ws.Range("AH8").Value2 = 'Blåbärsbarrskog'

# Row 9
$ws.Range("A9").Value2 = 112037208
$ws.Range("Q9").Value2 = 516097.2615754164
$ws.Range("R9").Value2 = 7184258.515744804
$ws.Range("Z9").Value2 = '11:44'
$ws.Range("AB9").Value2 = '11:44'
$ws.Range("AH9").Value2 = 'Blåbärsbarrskog'
$ws.Range("AM9").Value2 = 'Stående död trädstam/högstubbe'
$ws.Range("AO9").Value2 = 'Standing dead tree/snags'

# Row 10
$ws.Range("A10").Value2 = 112038473
$ws.Range("B10").Value2 = 89686
$ws.Range("E10").Value2 = 658
$ws.Range("F10").Value2 = 'Rosenticka'
$ws.Range("G10").Value2 = 'Rhodofomes roseus'
$ws.Range("H10").Value2 = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("I10").Value2 = '4'
$ws.Range("J10").Value2 = 'fruktkroppar'
$ws.Range("Q10").Value2 = 516057.2181607572
$ws.Range("R10").Value2 = 7184319.723381012
$ws.Range("Z10").Value2 = '13:34'
$ws.Range("AB10").Value2 = '13:34'
$ws.Range("AH10").Value2 = 'Blåbärsgranskog'

# Row 11
$ws.Range("A11").Value2 = 112038529
$ws.Range("Q11").Value2 = 515871.5299412137
$ws.Range("R11").Value2 = 7184628.386151251
$ws.Range("Z11").Value2 = '14:07'
$ws.Range("AB11").Value2 = '14:07'
$ws.Range("AH11").Value2 = 'Gransumpskog'
$ws.Range("AM11").ClearContents()
$ws.Range("AO11").ClearContents()

# Row 12
$ws.Range("A12").Value2 = 112037684
$ws.Range("B12").Value2 = 77515
$ws.Range("E12").Value2 = 6425
$ws.Range("F12").Value2 = 'Garnlav'
$ws.Range("G12").Value2 = 'Alectoria sarmentosa'
$ws.Range("H12").Value2 = '(Ach.) Ach.'
$ws.Range("Z12").Value2 = '12:08'
$ws.Range("AB12").Value2 = '12:08'
$ws.Range("AH12").Value2 = 'Blåbärsbarrskog'
$ws.Range("AM12").Value2 = 'Stående död trädstam/högstubbe'
$ws.Range("AO12").Value2 = 'Standing dead tree/snags # Picea abies'

# Row 13
$ws.Range("A13").Value2 = 112038436
$ws.Range("B13").Value2 = 89401
$ws.Range("E13").Value2 = 1108
$ws.Range("F13").Value2 = 'Harticka'
$ws.Range("G13").Value2 = 'Pelloporus leporinus'
$ws.Range("H13").Value2 = '(Fr.) Krieglst.'
$ws.Range("Q13").Value2 = 515951.3091604927
$ws.Range("R13").Value2 = 7184319.58691278
$ws.Range("Z13").Value2 = '13:28'
$ws.Range("AB13").Value2 = '13:28'
$ws.Range("AM13").Value2 = 'Stående död trädstam/högstubbe'
$ws.Range("AO13").Value2 = 'Standing dead tree/snags # Picea abies'
